$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class "1")
$ws.Range("B2").Value = 0.8130021913805697
$ws.Range("C2").Value = 0.9570077386070507
$ws.Range("D2").Value = 0.8791469194312795
$ws.Range("E2").Value = 1163

# Row 3 (class "2")
$ws.Range("B3").Value = 0.9528301886792453
$ws.Range("C3").Value = 0.9409937888198758
$ws.Range("D3").Value = 0.9468750000000001
$ws.Range("E3").Value = 644

# Row 4 (class "3")
$ws.Range("B4").Value = 0.8489116517285531
$ws.Range("C4").Value = 0.854381443298969
$ws.Range("D4").Value = 0.8516377649325625
$ws.Range("E4").Value = 776

# Row 5 (class "4")
$ws.Range("B5").Value = 0.935064935064935
$ws.Range("C5").Value = 0.4033613445378151
$ws.Range("D5").Value = 0.5636007827788649
$ws.Range("E5").Value = 357

# Row 6 (accuracy)
$ws.Range("B6").Value = 0.8591836734693877
$ws.Range("C6").Value = 0.8591836734693877
$ws.Range("D6").Value = 0.8591836734693877
$ws.Range("E6").Value = 0.8591836734693877

# Row 7 (macro avg)
$ws.Range("B7").Value = 0.8874522417133258
$ws.Range("C7").Value = 0.7889360788159276
$ws.Range("D7").Value = 0.8103151167856768

# Row 8 (weighted avg)
$ws.Range("B8").Value = 0.867931229130808
$ws.Range("C8").Value = 0.8591836734693877
$ws.Range("D8").Value = 0.8484053579381978
